# Auto-generated update of FFXIV leve market-price data (H:N columns)
# Mirrors a scheduled-runner refresh of currentAveragePrice / LevePrice / LeveProfit
# columns across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 62
$ws.Range("I62").Value = 1238.75
$ws.Range("J62").Value = 2300
$ws.Range("K62").Value = 1238.75
$ws.Range("L62").Value = 2300
$ws.Range("M62").Value = -614.75
$ws.Range("N62").Value = -3548
# ALC row 65
$ws.Range("I65").Value = 1238.75
$ws.Range("J65").Value = 2300
$ws.Range("K65").Value = 6193.75
$ws.Range("L65").Value = 11500
$ws.Range("M65").Value = -3073.75
$ws.Range("N65").Value = -17740
# ALC row 129
$ws.Range("H129").Value = 4611.7036
$ws.Range("I129").Value = 10484.7
$ws.Range("J129").Value = 1157
$ws.Range("K129").Value = 31454.1
$ws.Range("L129").Value = 3471
$ws.Range("M129").Value = -26454.1
$ws.Range("N129").Value = -13471

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 336493.34
$ws.Range("I2").Value = 5980
$ws.Range("J2").Value = 501750
$ws.Range("K2").Value = 5980
$ws.Range("L2").Value = 501750
$ws.Range("M2").Value = -5867
$ws.Range("N2").Value = -501976
# ARM row 61
$ws.Range("H61").Value = 1282.5625
$ws.Range("I61").Value = 1168.0667
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1168.0667
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -956.0667000000001
$ws.Range("N61").Value = -3424
# ARM row 116
$ws.Range("H116").Value = 336493.34
$ws.Range("I116").Value = 5980
$ws.Range("J116").Value = 501750
$ws.Range("K116").Value = 5980
$ws.Range("L116").Value = 501750
$ws.Range("M116").Value = -3686
$ws.Range("N116").Value = -506338
# ARM row 122
$ws.Range("H122").Value = 4550
$ws.Range("I122").Value = 4057.1428
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 12171.4284
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -9721.428400000001
$ws.Range("N122").Value = -28900
# ARM row 136
$ws.Range("H136").Value = 1282.5625
$ws.Range("I136").Value = 1168.0667
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3504.2001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -954.2001
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 336493.34
$ws.Range("I3").Value = 5980
$ws.Range("J3").Value = 501750
$ws.Range("K3").Value = 5980
$ws.Range("L3").Value = 501750
$ws.Range("M3").Value = -5866
$ws.Range("N3").Value = -501978
# BSM row 33
$ws.Range("H33").Value = 433.33334
$ws.Range("I33").Value = 433.33334
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 433.33334
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -97.33334000000002
$ws.Range("N33").ClearContents()
# BSM row 107
$ws.Range("H107").Value = 58824670
$ws.Range("I107").Value = 83334420
$ws.Range("J107").Value = 1279.8
$ws.Range("K107").Value = 83334420
$ws.Range("L107").Value = 1279.8
$ws.Range("M107").Value = -83332500
$ws.Range("N107").Value = -5119.8
# BSM row 122
$ws.Range("H122").Value = 40000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 40000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 40000
$ws.Range("N122").Value = -49800
# BSM row 134
$ws.Range("H134").Value = 2011.0186
$ws.Range("I134").Value = 1955.8914
$ws.Range("J134").Value = 2328
$ws.Range("K134").Value = 5867.674199999999
$ws.Range("L134").Value = 6984
$ws.Range("M134").Value = -3332.674199999999
$ws.Range("N134").Value = -12054

$ws = $wb.Worksheets.Item("CRP")
# CRP row 17
$ws.Range("H17").Value = 500
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -326
# CRP row 58
$ws.Range("H58").Value = 5327.1284
$ws.Range("I58").Value = 1452.25
$ws.Range("J58").Value = 15190.454
$ws.Range("K58").Value = 1452.25
$ws.Range("L58").Value = 15190.454
$ws.Range("M58").Value = -1249.25
$ws.Range("N58").Value = -15596.454
# CRP row 107
$ws.Range("H107").Value = 412.68
$ws.Range("I107").Value = 474.94116
$ws.Range("J107").Value = 280.375
$ws.Range("K107").Value = 474.94116
$ws.Range("L107").Value = 280.375
$ws.Range("M107").Value = 1445.05884
$ws.Range("N107").Value = -4120.375
# CRP row 132
$ws.Range("H132").Value = 2615.3901
$ws.Range("I132").Value = 3095.75
$ws.Range("J132").Value = 2307.96
$ws.Range("K132").Value = 9287.25
$ws.Range("L132").Value = 6923.88
$ws.Range("M132").Value = -6757.25
$ws.Range("N132").Value = -11983.88
# CRP row 136
$ws.Range("H136").Value = 5327.1284
$ws.Range("I136").Value = 1452.25
$ws.Range("J136").Value = 15190.454
$ws.Range("K136").Value = 4356.75
$ws.Range("L136").Value = 45571.362
$ws.Range("M136").Value = -1806.75
$ws.Range("N136").Value = -50671.362

$ws = $wb.Worksheets.Item("CUL")
# CUL row 125
$ws.Range("H125").Value = 1465
$ws.Range("I125").Value = 930
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 2790
$ws.Range("L125").Value = 6000
$ws.Range("M125").Value = 2130
$ws.Range("N125").Value = -15840

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70
$ws.Range("H70").Value = 60716.89
$ws.Range("I70").Value = 95149.77
$ws.Range("J70").Value = 6608.0713
$ws.Range("K70").Value = 95149.77
$ws.Range("L70").Value = 6608.0713
$ws.Range("M70").Value = -94879.77
$ws.Range("N70").Value = -7148.0713
# GSM row 73
$ws.Range("H73").Value = 60716.89
$ws.Range("I73").Value = 95149.77
$ws.Range("J73").Value = 6608.0713
$ws.Range("K73").Value = 95149.77
$ws.Range("L73").Value = 6608.0713
$ws.Range("M73").Value = -94213.77
$ws.Range("N73").Value = -8480.0713
# GSM row 107
$ws.Range("H107").Value = 449
$ws.Range("I107").Value = 368.125
$ws.Range("J107").Value = 772.5
$ws.Range("K107").Value = 368.125
$ws.Range("L107").Value = 772.5
$ws.Range("M107").Value = 1551.875
$ws.Range("N107").Value = -4612.5
# GSM row 120
$ws.Range("H120").Value = 35454.6
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 35454.6
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 35454.6
$ws.Range("N120").Value = -45130.6
# GSM row 122
$ws.Range("H122").Value = 2874.875
$ws.Range("I122").Value = 5750
$ws.Range("J122").Value = 1916.5
$ws.Range("K122").Value = 17250
$ws.Range("L122").Value = 5749.5
$ws.Range("M122").Value = -14800
$ws.Range("N122").Value = -10649.5

$ws = $wb.Worksheets.Item("LTW")
# LTW row 100
$ws.Range("H100").Value = 1394.1111
$ws.Range("I100").Value = 1325
$ws.Range("J100").Value = 1413.8572
$ws.Range("K100").Value = 1325
$ws.Range("L100").Value = 1413.8572
$ws.Range("M100").Value = -784
$ws.Range("N100").Value = -2495.8572

$ws = $wb.Worksheets.Item("WVR")
# WVR row 12
$ws.Range("H12").Value = 3050
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 3050
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3050
$ws.Range("N12").Value = -3334
# WVR row 101
$ws.Range("H101").Value = 12166.667
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 12166.667
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 12166.667
$ws.Range("N101").Value = -18656.667

